# Apply the edits described by the commit "first version of LCC calculation":
#  - Remove the "thermal bridge add on" column (column W) from Sheet1, together
#    with its header/unit/description/range cells. Deleting the entire column
#    shifts every later column left by one and Excel automatically garbage
#    collects the now-unused shared strings (4 strings removed from sst).
#  - Column O ("heat distribution") is narrowed from its old width to 20.
#  - Row 5 (the "required"/range notes row) is shortened from 75 to 60 points.
#  - The active selection moves from AG5 to X7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete entire column W (23rd column) -- this was the "thermal bridge add on"
# input column (header, unit "%", description and allowed-range helper text).
$ws.Columns.Item(23).Delete()

# Column O ("heat distribution") width change 64.14 -> 20
$ws.Columns.Item(15).ColumnWidth = 19.166666666666668

# Row 5 height change 75 -> 60
$ws.Rows.Item(5).RowHeight = 60

# Update the current selection to X7
$ws.Range("X7").Select() | Out-Null
